# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '26.859.52'
$ws.Cells.Item(2, 5).Value = '  +1.47%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.841.49'
$ws.Cells.Item(3, 5).Value = '  +1.65%  '
$ws.Cells.Item(4, 5).Value = '  +0.51%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '309.06'
$ws.Cells.Item(5, 5).Value = '  +1.13%  '
$ws.Cells.Item(6, 5).Value = '  +0.43%  '
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.4667'
$ws.Cells.Item(7, 5).Value = '  +3.72%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3653'
$ws.Cells.Item(8, 5).Value = '  +1.91%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.07136'
$ws.Cells.Item(9, 5).Value = '  +1.02%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.9136'
$ws.Cells.Item(10, 5).Value = '  +2.77%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07700'
$ws.Cells.Item(11, 5).Value = '  -0.96%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '19.52'
$ws.Cells.Item(12, 5).Value = '  +0.97%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.865.08'
$ws.Cells.Item(13, 5).Value = '  +3.12%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '5.280'
$ws.Cells.Item(14, 5).Value = '  +0.26%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '6.390'
$ws.Cells.Item(15, 5).Value = '  +1.34%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '87.99'
$ws.Cells.Item(16, 5).Value = '  +3.53%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '1.009'
$ws.Cells.Item(17, 5).Value = '  +0.42%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '0.000008668'
$ws.Cells.Item(18, 5).Value = '  +1.72%  '
$ws.Cells.Item(19, 5).Value = '  +0.45%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '26.887.22'
$ws.Cells.Item(20, 5).Value = '  +1.43%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '14.38'
$ws.Cells.Item(21, 5).Value = '  +1.69%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '5.018'
$ws.Cells.Item(22, 5).Value = '  +1.32%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '10.65'
$ws.Cells.Item(23, 5).Value = '  +1.45%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '1.938'
$ws.Cells.Item(24, 5).Value = '  -1.06%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '152.42'
$ws.Cells.Item(25, 5).Value = '  +1.21%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '18.24'
$ws.Cells.Item(26, 5).Value = '  +2.60%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '2.025'
$ws.Cells.Item(27, 5).Value = '  -1.31%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '114.06'
$ws.Cells.Item(28, 5).Value = '  +1.61%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '4.882'
$ws.Cells.Item(29, 5).Value = '  +1.05%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '0.08866'
$ws.Cells.Item(30, 5).Value = '  +2.15%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '3.206'
$ws.Cells.Item(31, 5).Value = '  +2.17%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.7458'
$ws.Cells.Item(32, 5).Value = '  +0.32%  '
$ws.Cells.Item(33, 2).Value = 'ARBITRUM'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.169'
$ws.Cells.Item(33, 5).Value = '  +5.58%  '
$ws.Cells.Item(34, 2).Value = 'RenderToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '2.783'
$ws.Cells.Item(34, 5).Value = '  +1.83%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '4.459'
$ws.Cells.Item(35, 5).Value = '  +0.68%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.083'
$ws.Cells.Item(36, 5).Value = '  +1.47%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.01943'
$ws.Cells.Item(37, 5).Value = '  +0.94%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '2.967'
$ws.Cells.Item(38, 5).Value = '  +2.52%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.05177'
$ws.Cells.Item(39, 5).Value = '  +1.67%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.5189'
$ws.Cells.Item(40, 5).Value = '  +1.99%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.903'
$ws.Cells.Item(41, 5).Value = '  +2.16%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.1513'
$ws.Cells.Item(42, 5).Value = '  +0.64%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '8.126'
$ws.Cells.Item(43, 5).Value = '  +1.02%  '
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '10.52'
$ws.Cells.Item(44, 5).Value = '  +5.60%  '
$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.4684'
$ws.Cells.Item(45, 5).Value = '  -0.64%  '
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.008'
$ws.Cells.Item(46, 5).Value = '  +0.52%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '100.38'
$ws.Cells.Item(47, 5).Value = '  +0.29%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '1.602'
$ws.Cells.Item(48, 5).Value = '  +1.71%  '
$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.06032'
$ws.Cells.Item(49, 5).Value = '  +0.83%  '
$ws.Cells.Item(50, 2).Value = 'Aave'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '64.59'
$ws.Cells.Item(50, 5).Value = '  +1.56%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.8818'
$ws.Cells.Item(51, 5).Value = '  +4.61%  '
